# Re-runs the "all algorithms" training sweep: the workbook gains a
# second sheet ("Accuracy Scores") summarising each model's accuracy, the
# detailed-report sheet is renamed and gains two new models (KNN, MLP),
# and every metric gets refreshed with the latest run's numbers.

$wb = $excel.ActiveWorkbook

# --- Sheets: rename Sheet1, add the new summary sheet right after it ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Detailed Reports"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Accuracy Scores"

# --- Sheet 1 "Detailed Reports": per-class precision/recall/f1 report ---
# per model (rows), sklearn classification_report-style dict literals.
$detailed = [ordered]@{
    "B1" = "0"; "C1" = "2"; "D1" = "3"; "E1" = "6"
    "F1" = "accuracy"; "G1" = "macro avg"; "H1" = "weighted avg"
    "I1" = "5"; "J1" = "loss"

    "A2" = "Logistic Regression"
    "B2" = "{'precision': 0.7777777777777778, 'recall': 0.7241379310344828, 'f1-score': 0.75, 'support': 29.0}"
    "C2" = "{'precision': 0.6363636363636364, 'recall': 0.7777777777777778, 'f1-score': 0.7, 'support': 18.0}"
    "D2" = "{'precision': 0.0, 'recall': 0.0, 'f1-score': 0.0, 'support': 2.0}"
    "E2" = "{'precision': 0.0, 'recall': 0.0, 'f1-score': 0.0, 'support': 1.0}"
    "F2" = 0.7
    "G2" = "{'precision': 0.35353535353535354, 'recall': 0.37547892720306514, 'f1-score': 0.3625, 'support': 50.0}"
    "H2" = "{'precision': 0.6802020202020203, 'recall': 0.7, 'f1-score': 0.687, 'support': 50.0}"
    "I2" = ""
    "J2" = ""

    "A3" = "KNN"
    "B3" = "{'precision': 0.7391304347826086, 'recall': 0.5862068965517241, 'f1-score': 0.6538461538461539, 'support': 29.0}"
    "C3" = "{'precision': 0.5789473684210527, 'recall': 0.6111111111111112, 'f1-score': 0.5945945945945946, 'support': 18.0}"
    "D3" = "{'precision': 0.0, 'recall': 0.0, 'f1-score': 0.0, 'support': 2.0}"
    "E3" = "{'precision': 0.3333333333333333, 'recall': 1.0, 'f1-score': 0.5, 'support': 1.0}"
    "F3" = 0.58
    "G3" = "{'precision': 0.3302822273073989, 'recall': 0.4394636015325671, 'f1-score': 0.3496881496881497, 'support': 50.0}"
    "H3" = "{'precision': 0.6437833714721587, 'recall': 0.58, 'f1-score': 0.6032848232848234, 'support': 50.0}"
    "I3" = "{'precision': 0.0, 'recall': 0.0, 'f1-score': 0.0, 'support': 0.0}"
    "J3" = ""

    "A4" = "SVM"
    "B4" = "{'precision': 0.7777777777777778, 'recall': 0.7241379310344828, 'f1-score': 0.75, 'support': 29.0}"
    "C4" = "{'precision': 0.5714285714285714, 'recall': 0.6666666666666666, 'f1-score': 0.6153846153846154, 'support': 18.0}"
    "D4" = "{'precision': 0.0, 'recall': 0.0, 'f1-score': 0.0, 'support': 2.0}"
    "E4" = "{'precision': 0.5, 'recall': 1.0, 'f1-score': 0.6666666666666666, 'support': 1.0}"
    "F4" = 0.68
    "G4" = "{'precision': 0.4623015873015873, 'recall': 0.5977011494252873, 'f1-score': 0.5080128205128205, 'support': 50.0}"
    "H4" = "{'precision': 0.6668253968253969, 'recall': 0.68, 'f1-score': 0.6698717948717949, 'support': 50.0}"
    "I4" = ""
    "J4" = ""

    "A5" = "Random Forest"
    "B5" = "{'precision': 0.8214285714285714, 'recall': 0.7931034482758621, 'f1-score': 0.8070175438596491, 'support': 29.0}"
    "C5" = "{'precision': 0.7, 'recall': 0.7777777777777778, 'f1-score': 0.7368421052631579, 'support': 18.0}"
    "D5" = "{'precision': 1.0, 'recall': 0.5, 'f1-score': 0.6666666666666666, 'support': 2.0}"
    "E5" = "{'precision': 1.0, 'recall': 1.0, 'f1-score': 1.0, 'support': 1.0}"
    "F5" = 0.78
    "G5" = "{'precision': 0.8803571428571428, 'recall': 0.76772030651341, 'f1-score': 0.8026315789473684, 'support': 50.0}"
    "H5" = "{'precision': 0.7884285714285714, 'recall': 0.78, 'f1-score': 0.78, 'support': 50.0}"
    "I5" = ""
    "J5" = ""

    "A6" = "MLP"
    "B6" = "{'precision': 0.84, 'recall': 0.7241379310344828, 'f1-score': 0.7777777777777778, 'support': 29.0}"
    "C6" = "{'precision': 0.6842105263157895, 'recall': 0.7222222222222222, 'f1-score': 0.7027027027027027, 'support': 18.0}"
    "D6" = "{'precision': 0.6666666666666666, 'recall': 1.0, 'f1-score': 0.8, 'support': 2.0}"
    "E6" = "{'precision': 0.5, 'recall': 1.0, 'f1-score': 0.6666666666666666, 'support': 1.0}"
    "F6" = 0.74
    "G6" = "{'precision': 0.5381754385964912, 'recall': 0.689272030651341, 'f1-score': 0.5894294294294294, 'support': 50.0}"
    "H6" = "{'precision': 0.770182456140351, 'recall': 0.74, 'f1-score': 0.7494174174174174, 'support': 50.0}"
    "I6" = "{'precision': 0.0, 'recall': 0.0, 'f1-score': 0.0, 'support': 0.0}"
    "J6" = ""

    "A7" = "Neural Network"
    "B7" = ""
    "C7" = ""
    "D7" = ""
    "E7" = ""
    "F7" = 0.6800000071525574
    "G7" = ""
    "H7" = ""
    "I7" = ""
    "J7" = 0.8243193030357361
}

foreach ($addr in $detailed.Keys) {
    $ws1.Range($addr).Value = $detailed[$addr]
}

# --- Sheet 2 "Accuracy Scores": one row per model, overall accuracy only ---
$accuracy = [ordered]@{
    "A1" = "Model"; "B1" = "Accuracy"
    "A2" = "Logistic Regression"; "B2" = 0.7
    "A3" = "KNN";                 "B3" = 0.58
    "A4" = "SVM";                 "B4" = 0.68
    "A5" = "Random Forest";       "B5" = 0.78
    "A6" = "MLP";                 "B6" = 0.74
    "A7" = "Neural Network";      "B7" = 0.6800000071525574
}

foreach ($addr in $accuracy.Keys) {
    $ws2.Range($addr).Value = $accuracy[$addr]
}

# --- Formatting: the header row (sheet1 row 1) and both sheets' leftmost
# "label" column already carry a bold/boxed/centered style (style index 1
# in the original file, reused by the model names already on sheet1).
# Copy that format onto every header / label cell so newly-written ones
# (J1, A6, A7, sheet2's whole A/row1) pick it up too, instead of Excel
# minting a brand-new near-duplicate style.
$fmtSource = $ws1.Range("A2")
$fmtSource.Copy()

$ws1.Range("B1:J1").PasteSpecial(-4122)   # xlPasteFormats
$ws1.Range("A2:A7").PasteSpecial(-4122)
$ws2.Range("A1:B1").PasteSpecial(-4122)
$ws2.Range("A2:A7").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Leave "Detailed Reports" as the active/selected tab, matching its
# original position as sheet 1.
$ws1.Activate()
